$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

# 1. Heading3 title: add " I" (scope to the Heading 3 paragraph to avoid
#    touching the unrelated "Atmospheric pollution:" bibliography entry)
$heading3 = $d.Paragraphs(2).Range
$heading3.Find.Execute("Atmospheric pollution", $true, $true, $false, $false, `
                        $false, $true, 1, $false, "Atmospheric pollution I", 2)

# 2. Activation date
Replace-Text "Ativação: 01/01/2022" "Ativação: 01/01/2025"

# 3. Objectives (English, italic)
Replace-Text "Enable the student to identify the main pollutants from the atmosphere, to evaluate their influence on climate change and make the management and monitoring of air quality in large Brazilian cities." "Enable the student to identify the main pollutants in the atmosphere and their effects on man and the environment"

# 4. Docente responsável
Replace-Text "7455355 - Robson da Silva Rocha" "4893449 - Débora Souza Alvim"

# 5. Programa resumido (Portuguese)
Replace-Text "Características e composição da atmosfera. Origem, movimentação e destino dos poluentes. Histórico da poluição do ar. Principais poluentes atmosféricos e padrões da qualidade do ar.  Dispersão de poluentes na atmosfera. Modelos matemáticos do transporte de poluentes atmosféricos. Qualidade do ar no interior de edifícios. Controle de poluição de fontes fixas e móveis. Otimização exergoeconômica/ambiental." "1)Características e composição da atmosfera. Origem, movimentação e destino dos poluentes. Histórico da poluição do ar. Principais poluentes atmosféricos e padrões da qualidade do ar.  O efeito estufa. Dispersão de poluentes na atmosfera. Modelos matemáticos do transporte de poluentes atmosféricos.  Controle da poluição atmosférica de fontes fixas e móveis."

# 6. Programa resumido (English, italic) - replaced with combined PT+EN text
Replace-Text "Characteristics and composition of the atmosphere. Origin, movement and fate of pollutants. History of air pollution. Major air pollutants and air quality standards. Dispersion of pollutants in the atmosphere. Mathematical models of transport of air pollutants. Air quality inside buildings. Control pollution of fixed and mobile sources. Exergoeconomic / environmental optimization" "1)Características e composição da atmosfera. Origem, movimentação e destino dos poluentes. Histórico da poluição do ar. Principais poluentes atmosféricos e padrões da qualidade do ar.  O efeito estufa. Dispersão de poluentes na atmosfera. Modelos matemáticos do transporte de poluentes atmosféricos.  Controle da poluição atmosférica de fontes fixas e móveis.1)Characteristics and composition of the atmosphere. Origin, movement and fate of pollutants. History of air pollution. Major air pollutants and air quality standards. The greenhouse effect. Dispersion of pollutants in the atmosphere. Mathematical models of transport of air pollutants. Air Pollution control from fixed and mobile sources."

# 7. Programa (Portuguese)
Replace-Text "1) Caracterização da atmosfera e seus poluentes. 2) Padrões da qualidade do ar. 3) Dispersão de poluentes na atmosfera. 4) Modelagem matemática do transporte de poluentes.5) Qualidade do ar no interior de edifícios.6) Controle de poluição de fontes fixas e móveis.7) Otimização exergoeconômica/ambiental." "Caracterização da atmosfera e seus poluentes. 2) Padrões da qualidade do ar. 3) Dispersão de poluentes na atmosfera. 4) O efeito estufa. 5) Modelagem matemática do transporte de poluentes. 6) Controle da poluição atmosférica de fontes fixas e móveisA disciplina pode contar com viagens didáticas para complementação do conteúdo da disciplina."

# 8. Programa (English, italic)
Replace-Text "Characterization of the atmosphere and its pollutants.Air quality standards.Dispersion of pollutants in the atmosphere.Mathematical modeling of pollutant transport.Air quality inside buildings.Control pollution of fixed and mobile sourcesExergoeconomic / environmental optimization." "Characterization of the atmosphere and its pollutants. 2) Air quality standards. 3) Dispersion of pollutants in the atmosphere. 4) The greenhouse effect. 5) Mathematical modeling of pollutant transport. 6) Air Pollution control from fixed and mobile sources. The discipline may have didactic trips to complement the content of the discipline."

# 9. Bibliografia
Replace-Text "Bibliografia básica:GUNTER, F.; Introdução aos problemas da poluição ambiental. 1 ed. São Paulo: Editora EPU, 2008.LENZI, E. F.; FAVERO, L.O.B. Introdução à química da atmosfera  Ciência, vida e sobrevivência. 1ª. ed. Rio de Janeiro: Editora LCT, 465p. 2009.SPIRO, T. G.; STIGLIANI, E. W. M. Química ambiental. 2 ed. Sao Paulo: Pearson / Prentice Hall. 2008. 352p.VESILIND, P. A.; MORGAN, S. M., revisão técnica Carlos Alberto de Moya Figueira Netto; Lineu Belico dos Reis. Introdução à Engenharia Ambiental. Tradução da 2ª edição norte-americana. Editora Cengage Learning, São Paulo, 2015.Bibliografia complementar:JACOBSON, Mark Z. Atmospheric pollution: history, science, and regulation. Cambridge, Inglaterra: Cambridge University Press, c2002. xi, 399 p. Includes bibliographical references (p 355-370). ISBN 9780521010443.SEINFELD, J.H.;MANAHAN, S.E. Environmental chemistry. 9 th edition. Boca Raton, FL: CRC Press. 753p. 2010.PANDIS, S.N. Atmospheric Chemistry and Physics: From air pollution to climate change. John Wiley& Sons, 1998.SCHNELLE JR, Karl B; BROWN, Charles A. Air pollution control technology handbook. New York: CRC Press, 2001. 386 p. (Mechanical engineering handbook series). ISBN 9780849395888.VALLERO, Daniel A. Fundamentals of air pollution. 4 ed. Amsterdam: Elsevier, 2008. 942 p" "Bibliografia básica:Baird, C.; Cann, M. Química Ambiental. Porto Alegre: Bookman, 4.ed., 2011. 844p.GUNTER, F.; Introdução aos problemas da poluição ambiental. 1 ed. São Paulo: Editora EPU, 2008.LENZI, E. F.; FAVERO, L.O.B. Introdução à química da atmosfera  Ciência, vida e sobrevivência. 1ª. ed. Rio de Janeiro: Editora LCT, 465p. 2009.Rocha, Julio Cesar; Rosa, André Henrique; Cardoso, Arnaldo Alves. Introdução à química ambiental. 2. ed. Porto Alegre: Bookman, 2009.  03Seinfeld, J.H. e Pandis, S.P. Atmospheric Chemistry and Physics: from air pollution to climate change. New York, USA: John Wiley & Sons Inc., 2006.SPIRO, T. G.; STIGLIANI, E. W. M. Química ambiental. 2 ed. Sao Paulo: Pearson / Prentice Hall. 2008. 352p."
